$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.613.04"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.863.80"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.26"
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4621"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3908"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07888"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9683"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.26"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").Value = "1.865.93"
$ws.Range("E12").Value = "  +6.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.718"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.920"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06930"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.41"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001003"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.91"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D21").Value = "28.625.28"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.316"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.05"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.123"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").Value = "2.075.12"
$ws.Range("E25").Value = "  +3.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.99"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.29"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.749"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.991"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.03"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09341"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9350"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.312"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.336"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.345"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05811"
$ws.Range("E36").Value = "  -3.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02113"
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.155"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.929"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5648"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.917"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1775"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.68"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.185"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5304"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.139"
$ws.Range("E47").Value = "  -8.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.846"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.31"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.346"
$ws.Range("E51").Value = "  +0.83%  "